$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "Test Comment" column (K) descriptions added for this export
$ws.Range('K8').Value = 'General fault alarm as a result of loss of power'
$ws.Range('K9').Value = 'CYP holding doors open on site'
$ws.Range('K10').Value = 'Door forced open'
$ws.Range('K17').Value = 'Door forced open'
$ws.Range('K18').Value = 'CYP holding doors open on site'
$ws.Range('K38').Value = 'CYP holding doors open on site'
$ws.Range('K39').Value = 'PED being used for access'
$ws.Range('K47').Value = 'PED being used for access'
$ws.Range('K48').Value = 'Isolation of lift causes this alarm'
$ws.Range('K49').Value = 'Isolation of lift causes this alarm'
$ws.Range('K54').Value = 'Isolation of lift causes this alarm'
$ws.Range('K60').Value = 'Isolation of lift causes this alarm'
$ws.Range('K61').Value = 'Isolation of panel'
$ws.Range('K111').Value = 'PED being used for access'
$ws.Range('K113').Value = 'PED being used for access'
$ws.Range('K114').Value = 'Door forced open'
$ws.Range('K115').Value = 'CYP holding doors open on site'
$ws.Range('K118').Value = 'CYP holding doors open on site'
$ws.Range('K146').Value = 'Door forced open'
$ws.Range('K147').Value = 'Door forced open'
$ws.Range('K148').Value = 'CYP holding doors open on site'
$ws.Range('K149').Value = 'PSD isolated during a blockade'
$ws.Range('K168').Value = 'PSD opening'
$ws.Range('K169').Value = 'PSD isolated during a blockade'
$ws.Range('K170').Value = 'PSD isolated during a blockade'
$ws.Range('K171').Value = 'PSD isolated during a blockade'
$ws.Range('K172').Value = 'PSD isolated during a blockade'
$ws.Range('K173').Value = 'PSD opening'
$ws.Range('K185').Value = 'PSD isolated during a blockade'
$ws.Range('K186').Value = 'PSD isolated during a blockade'
$ws.Range('K187').Value = 'PSD isolated during a blockade'
$ws.Range('K188').Value = 'Real alarm to be looked at'
$ws.Range('K189').Value = 'CYP holding doors open on site'
$ws.Range('K190').Value = 'CYP holding doors open on site'
$ws.Range('K258').Value = 'Not currently connected to Sunshine'
$ws.Range('K259').Value = 'Not currently connected to Sunshine'
$ws.Range('K260').Value = 'Not currently connected to Sunshine'
$ws.Range('K261').Value = 'Not currently connected to Sunshine'
$ws.Range('K262').Value = 'Not currently connected to Sunshine'
$ws.Range('K277').Value = 'Not currently connected to Sunshine'
$ws.Range('K279').Value = 'Not currently connected to Sunshine'

# Size the new column to fit its content (matches bestFit column K width)
$ws.Columns.Item(11).ColumnWidth = 40.16666666666667

# Re-create the AutoFilter over the widened A1:K280 range, preserving the
# existing Location ("ARN") filter, and drop the stale sort state
$ws.AutoFilterMode = $false
$ws.Range('A1:K280').AutoFilter(10, @('ARN')) | Out-Null

# Keep the _FilterDatabase defined name in sync with the new filter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq 'Sheet1!_FilterDatabase') { $n.RefersTo = '=Sheet1!$A$1:$K$280' }
}

# Restore the active cell/selection to match the author's final position
$ws.Range('K49').Select() | Out-Null
